$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049429957228124
$ws.Range("D2").Value = 1.048819649035727
$ws.Range("E2").Value = 1.060514091990276
$ws.Range("F2").Value = 1.068508330146617
$ws.Range("I2").Value = 1.044828802889624
$ws.Range("J2").Value = 1.054468447793639
$ws.Range("K2").Value = 1.051578430221607
$ws.Range("L2").Value = 1.063240664831152
$ws.Range("M2").Value = 1.071213324540388
$ws.Range("N2").Value = 1.055965913938727

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050535264081196
$ws.Range("D3").Value = 1.04964321193144
$ws.Range("E3").Value = 1.061689931887884
$ws.Range("F3").Value = 1.069732729822327
$ws.Range("I3").Value = 1.045147836785625
$ws.Range("J3").Value = 1.055222268688761
$ws.Range("K3").Value = 1.052214085521321
$ws.Range("L3").Value = 1.064230032790796
$ws.Range("M3").Value = 1.07225270133949
$ws.Range("N3").Value = 1.056720805345985

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.051250233543303
$ws.Range("D4").Value = 1.05017589164926
$ws.Range("E4").Value = 1.062451483468579
$ws.Range("F4").Value = 1.070525608855001
$ws.Range("I4").Value = 1.045353009753511
$ws.Range("J4").Value = 1.055709231966445
$ws.Range("K4").Value = 1.052624527765913
$ws.Range("L4").Value = 1.064870345475302
$ws.Range("M4").Value = 1.072925286946877
$ws.Range("N4").Value = 1.057208460167378

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051550750663832
$ws.Range("D5").Value = 1.050399777394667
$ws.Range("E5").Value = 1.062771809213165
$ws.Range("F5").Value = 1.070859082259781
$ws.Range("I5").Value = 1.045438962173971
$ws.Range("J5").Value = 1.055913757983372
$ws.Range("K5").Value = 1.05279686978917
$ws.Range("L5").Value = 1.065139563661296
$ws.Range("M5").Value = 1.073208052094859
$ws.Range("N5").Value = 1.057413276634709

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051601205551561
$ws.Range("D6").Value = 1.050437365699502
$ws.Range("E6").Value = 1.062825603275509
$ws.Range("F6").Value = 1.070915082561345
$ws.Range("I6").Value = 1.04545337624094
$ws.Range("J6").Value = 1.055948087475097
$ws.Range("K6").Value = 1.052825794587927
$ws.Range("L6").Value = 1.065184768401958
$ws.Range("M6").Value = 1.073255530224909
$ws.Range("N6").Value = 1.05744765487825

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.051254249286806
$ws.Range("D7").Value = 1.050178883429061
$ws.Range("E7").Value = 1.062455763009993
$ws.Range("F7").Value = 1.070530064163992
$ws.Range("I7").Value = 1.045354159441596
$ws.Range("J7").Value = 1.055711965612318
$ws.Range("K7").Value = 1.052626831425103
$ws.Range("L7").Value = 1.064873942660237
$ws.Range("M7").Value = 1.072929065227766
$ws.Range("N7").Value = 1.057211197695342

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049803550690548
$ws.Range("D8").Value = 1.0490980216334
$ws.Range("E8").Value = 1.060911327184691
$ws.Range("F8").Value = 1.068921995860803
$ws.Range("I8").Value = 1.044936883546758
$ws.Range("J8").Value = 1.05472337289008
$ws.Range("K8").Value = 1.051793433009672
$ws.Range("L8").Value = 1.063575000808284
$ws.Range("M8").Value = 1.071564579103393
$ws.Range("N8").Value = 1.056221201058045

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047245369686502
$ws.Range("D9").Value = 1.047191715696341
$ws.Range("E9").Value = 1.058195194204527
$ws.Range("F9").Value = 1.066093010294766
$ws.Range("I9").Value = 1.044191908254254
$ws.Range("J9").Value = 1.052975130423456
$ws.Range("K9").Value = 1.050318211735398
$ws.Range("L9").Value = 1.061287028662255
$ws.Range("M9").Value = 1.069160447697016
$ws.Range("N9").Value = 1.054470475886593

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045538600449988
$ws.Range("D10").Value = 1.045919705993091
$ws.Range("E10").Value = 1.056387990483276
$ws.Range("F10").Value = 1.064210078580921
$ws.Range("I10").Value = 1.043688738051682
$ws.Range("J10").Value = 1.051805428864607
$ws.Range("K10").Value = 1.049330233032263
$ws.Range("L10").Value = 1.059762292479187
$ws.Range("M10").Value = 1.067557823840978
$ws.Range("N10").Value = 1.053299113217364

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044799224040663
$ws.Range("D11").Value = 1.045368639935172
$ws.Range("E11").Value = 1.055606281283435
$ws.Range("F11").Value = 1.063395460008315
$ws.Range("I11").Value = 1.043469310562002
$ws.Range("J11").Value = 1.051297930560629
$ws.Range("K11").Value = 1.0489013563484
$ws.Range("L11").Value = 1.059102191105598
$ws.Range("M11").Value = 1.066863889193958
$ws.Range("N11").Value = 1.052790894207583

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044524535172393
$ws.Range("D12").Value = 1.04516390753326
$ws.Range("E12").Value = 1.055316042005545
$ws.Range("F12").Value = 1.063092978772048
$ws.Range("I12").Value = 1.043387572025551
$ws.Range("J12").Value = 1.051109270690899
$ws.Range("K12").Value = 1.048741890436038
$ws.Range("L12").Value = 1.058857016974832
$ws.Range("M12").Value = 1.066606131672137
$ws.Range("N12").Value = 1.052601966419199

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.044583459221633
$ws.Range("D13").Value = 1.045207825229941
$ws.Range("E13").Value = 1.0553782937991
$ws.Range("F13").Value = 1.063157857321622
$ws.Range("I13").Value = 1.043405115785689
$ws.Range("J13").Value = 1.051149745781051
$ws.Range("K13").Value = 1.048776103754405
$ws.Range("L13").Value = 1.058909606909541
$ws.Range("M13").Value = 1.066661421521491
$ws.Range("N13").Value = 1.052642498988621

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044776519232571
$ws.Range("D14").Value = 1.045351717547794
$ws.Range("E14").Value = 1.055582287525516
$ws.Range("F14").Value = 1.063370454693242
$ws.Range("I14").Value = 1.04346255879207
$ws.Range("J14").Value = 1.051282338989922
$ws.Range("K14").Value = 1.048888178154038
$ws.Range("L14").Value = 1.059081924596486
$ws.Range("M14").Value = 1.066842582868848
$ws.Range("N14").Value = 1.052775280495058

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044895463005672
$ws.Range("D15").Value = 1.045440368790086
$ws.Range("E15").Value = 1.055707990967455
$ws.Range("F15").Value = 1.063501456742362
$ws.Range("I15").Value = 1.043497920387687
$ws.Range("J15").Value = 1.051364013826449
$ws.Range("K15").Value = 1.048957209408318
$ws.Range("L15").Value = 1.059188097420435
$ws.Range("M15").Value = 1.066954202414139
$ws.Range("N15").Value = 1.05285707131922

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045587663272776
$ws.Range("D16").Value = 1.045956272565706
$ws.Range("E16").Value = 1.056439887166407
$ws.Range("F16").Value = 1.064264156840128
$ws.Range("I16").Value = 1.043703268017677
$ws.Range("J16").Value = 1.051839088521165
$ws.Range("K16").Value = 1.049358673469635
$ws.Range("L16").Value = 1.059806103661648
$ws.Range("M16").Value = 1.067603878157522
$ws.Range("N16").Value = 1.053332820674495

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046021771814482
$ws.Range("D17").Value = 1.046279810917052
$ws.Range("E17").Value = 1.056899205550995
$ws.Range("F17").Value = 1.064742765841225
$ws.Range("I17").Value = 1.04383166157198
$ws.Range("J17").Value = 1.052136819484073
$ws.Range("K17").Value = 1.049610213039299
$ws.Range("L17").Value = 1.060193793899112
$ws.Range("M17").Value = 1.068011405033222
$ws.Range("N17").Value = 1.05363097444952

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046274947840953
$ws.Range("D18").Value = 1.04646849860917
$ws.Range("E18").Value = 1.057167197576601
$ws.Range("F18").Value = 1.065021998317295
$ws.Range("I18").Value = 1.043906401685938
$ws.Range("J18").Value = 1.052310383531901
$ws.Range("K18").Value = 1.049756828085831
$ws.Range("L18").Value = 1.060419938650592
$ws.Range("M18").Value = 1.068249109802733
$ws.Range("N18").Value = 1.053804784978204

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046361268912921
$ws.Range("D19").Value = 1.046532831705483
$ws.Range("E19").Value = 1.057258589431887
$ws.Range("F19").Value = 1.065117220939154
$ws.Range("I19").Value = 1.043931860755335
$ws.Range("J19").Value = 1.05236954788961
$ws.Range("K19").Value = 1.049806802447348
$ws.Range("L19").Value = 1.060497050165029
$ws.Range("M19").Value = 1.068330161229957
$ws.Range("N19").Value = 1.053864033356088

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045975199410454
$ws.Range("D20").Value = 1.046245101082866
$ws.Range("E20").Value = 1.05684991683253
$ws.Range("F20").Value = 1.064691408606337
$ws.Range("I20").Value = 1.043817901642326
$ws.Range("J20").Value = 1.052104885847341
$ws.Range("K20").Value = 1.049583235956308
$ws.Range("L20").Value = 1.060152197213667
$ws.Range("M20").Value = 1.06796768113605
$ws.Range("N20").Value = 1.053598995463362

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044719669305171
$ws.Range("D21").Value = 1.045309346006124
$ws.Range("E21").Value = 1.055522213090068
$ws.Range("F21").Value = 1.063307847195297
$ws.Range("I21").Value = 1.043445649703767
$ws.Range("J21").Value = 1.051243297807337
$ws.Range("K21").Value = 1.048855179512803
$ws.Range("L21").Value = 1.059031180874093
$ws.Range("M21").Value = 1.066789235356851
$ws.Range("N21").Value = 1.052736183869515

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043929968775535
$ws.Range("D22").Value = 1.04472075692084
$ws.Range("E22").Value = 1.054688140208445
$ws.Range("F22").Value = 1.062438550453123
$ws.Range("I22").Value = 1.043210249622603
$ws.Range("J22").Value = 1.050700701109234
$ws.Range("K22").Value = 1.048396483929691
$ws.Range("L22").Value = 1.058326450734751
$ws.Range("M22").Value = 1.066048303245875
$ws.Range("N22").Value = 1.052192816621863

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.044348632756349
$ws.Range("D23").Value = 1.0450328022224
$ws.Range("E23").Value = 1.05513023138166
$ws.Range("F23").Value = 1.062899324354446
$ws.Range("I23").Value = 1.043335167774567
$ws.Range("J23").Value = 1.050988425729183
$ws.Range("K23").Value = 1.048639736162087
$ws.Range("L23").Value = 1.058700032635954
$ws.Range("M23").Value = 1.066441085412681
$ws.Range("N23").Value = 1.052480949843777

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045996243574051
$ws.Range("D24").Value = 1.046260785045358
$ws.Range("E24").Value = 1.056872188037157
$ws.Range("F24").Value = 1.064714614519811
$ws.Range("I24").Value = 1.043824119624584
$ws.Range("J24").Value = 1.05211931558386
$ws.Range("K24").Value = 1.049595426059202
$ws.Range("L24").Value = 1.060170992929157
$ws.Range("M24").Value = 1.067987438078934
$ws.Range("N24").Value = 1.053613445691762

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04790694810524
$ws.Range("D25").Value = 1.047684741993373
$ws.Range("E25").Value = 1.058896749861693
$ws.Range("F25").Value = 1.066823827483695
$ws.Range("I25").Value = 1.044385650652811
$ws.Range("J25").Value = 1.053427832428457
$ws.Range("K25").Value = 1.050700383431613
$ws.Range("L25").Value = 1.061878418940674
$ws.Range("M25").Value = 1.069781946689026
$ws.Range("N25").Value = 1.054923820780367
